$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header for new column F ("scenario"), matching the bold/centered
# formatting already used by the other header cells (A1:E1)
$ws.Range("F1").Value = "scenario"
$ws.Range("F1").Font.Bold = $true
$ws.Range("F1").HorizontalAlignment = -4108

# Fill F2:F101 with the scenario label "S1" for every data row
$ws.Range("F2:F101").Value = "S1"
